$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1 (header cell)
$ws.Range("F1").Value = "Last status check on: 20.02.2022 23:45"

# Row 10 (EuroOil Opustena) got a fresh price-check run:
# the previous "new" price (B10) becomes the "old" price (C10) and
# a newly scraped price takes over B10 - i.e. the two values swap.
$ws.Range("B10").Value = 36.5
$ws.Range("C10").Value = 37.4

# Delta Cena flips sign (price dropped instead of rose) - keep it a
# plain text value like the original ("+0.9"), not a numeric cell.
$ws.Range("D10").Value = "'-0.9"
$ws.Range("D10").Style = "Normal"

# Status-check timestamp text for this row is updated
$ws.Range("E10").Value = "2022-02-20 23:47:14"
